{"js": "// Diff summary: the document gains two new paragraphs at the very end\n// (a new \"19.05.2021\" log date, followed by its bullet entry about the\n// new weapon class). Every other hunk in the diff is a no-op at the text\n// level (Word merely re-merged previously split runs / dropped stale\n// <w:proofErr> spell-check markers while resaving) so there is nothing\n// else to replicate here.\nconst body = context.document.body;\n\n// New date heading paragraph, appended after the last paragraph in the doc.\nbody.insertParagraph(\"19.05.2021\", Word.InsertLocation.end);\n\n// New bullet paragraph describing the weapon class work.\nbody.insertParagraph(\n  \"-Waffenklasse erstellt, subKlassen f\u00fcr verschiedene Waffentypen erstellt. Aktuell ausger\u00fcstete Waffe wird auf dem Bildschirm angezeigt\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# Diff summary: the document gains two new paragraphs at the very end\n# (a new \"19.05.2021\" log date, followed by its bullet entry about the\n# new weapon class). Every other hunk in the diff is a no-op at the text\n# level (Word merely re-merged previously split runs / dropped stale\n# proofing marks while resaving) so there is nothing else to replicate\n# here.\n$d = $word.ActiveDocument\n\n# Append the new date heading paragraph after the current last paragraph.\n$tail = $d.Paragraphs.Last.Range\n$tail.InsertParagraphAfter()\n$datePara = $d.Paragraphs.Last\n$datePara.Range.Text = \"19.05.2021\"\n\n# Append the new bullet paragraph describing the weapon class work.\n$tail2 = $d.Paragraphs.Last.Range\n$tail2.InsertParagraphAfter()\n$bulletPara = $d.Paragraphs.Last\n$bulletPara.Range.Text = \"-Waffenklasse erstellt, subKlassen f\u00fcr verschiedene Waffentypen erstellt. Aktuell ausger\u00fcstete Waffe wird auf dem Bildschirm angezeigt\"\n"}
